# ajustes sanity semilla 6 en clases de portabilidad prepago y postpago
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semilla 6")

# Update numeric/text reference values for prepago/postpago classes
$ws.Range("H9").Value = "662496115"
$ws.Range("E10").Value = "3046008600"

# Update portability service URL (prepago) from .76 to .74 host
$ws.Range("I2").Value = "http://10.69.60.74:8080/PortabilidadServiceEAR-HPNPCommunicationsDelegateEJB/NPCRMWSImpl"

$ws.Range("B14").Value = "662496115"

# Restore the window/selection state captured in the workbook
$ws.Range("I2").Select()

$win = $excel.ActiveWindow
$win.ScrollColumn = 3
